# Change the table style of the table on slide 16 (the "PLENARY- COMPLETE
# THE MISSING GAPS" slide) from the custom "Table_0" style to the built-in
# table style identified by {2C3633B0-DAE9-443D-903A-D34F20F63772}, exactly
# as PowerPoint records it when a different style is picked from the
# Table Design > Table Styles gallery.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)

# Locate the shape that holds the table (it is the 3rd shape on the slide).
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $tableShape = $shape
        break
    }
}

$table = $tableShape.Table
$table.ApplyStyle("{2C3633B0-DAE9-443D-903A-D34F20F63772}", $true)
